# Generate Report for Handoff
# Updates the localization-status report: the "In Translation" status moves to
# "Ready for handoff" and the handoff timestamps are refreshed, on all three
# sheets (Overview, zh-cn, de-de). Also widens the Status/zh-cn/de-de columns
# to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Handoff timestamps refreshed to the new generate time ------------------
$overview.Range("G2").Value = "2016-08-26 19:01:43"
$zhcn.Range("H2").Value     = "2016-08-26 19:01:38"
$dede.Range("H2").Value     = "2016-08-26 19:01:43"

# --- Widen the Status / zh-cn / de-de columns to fit "Ready for handoff" ----
$overview.Range("E1:F1").ColumnWidth = 16.33
$zhcn.Range("C1").ColumnWidth = 16.33
$dede.Range("C1").ColumnWidth = 16.33
